$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Language dropdown rows (10 = en, 11 = fr, 12 = ar): show each language's
# own name in every translation column instead of mixed language names ---

# Row 10 (English): en / English / English / English
$ws.Range("E10").Value = "English"
$ws.Range("F10").Value = "English"

# Row 11 (French): fr / Français / Français / Français
$ws.Range("D11").Value = "Français"
$ws.Range("F11").Value = "Français"

# Row 12 (Arabic): ar / العربية / العربية / العربية
$ws.Range("D12").Value = "العربية"
$ws.Range("E12").Value = "العربية"

# Match the cell borders to the rest of the table: columns E/F normally carry
# a full thin border (same style as D10 already has) while D11/D12 now need
# the same border treatment already used on E11/E12.
$ws.Range("D10").Copy()
$ws.Range("E10:F10").PasteSpecial(-4122)

$ws.Range("E11").Copy()
$ws.Range("D11").PasteSpecial(-4122)

$ws.Range("E12").Copy()
$ws.Range("D12").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Move the active selection from the old spot to the new focal cell ---
$ws.Range("F14").Select()
